$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 20:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1112341
$ws.Range("C4").Value = 17318
$ws.Range("D4").Value = 157811
$ws.Range("E4").Value = 889616
$ws.Range("G4").Value = 1058
$ws.Range("H4").Value = 64914

# Row 8 - Francia
$ws.Range("D8").Value = 50212
$ws.Range("E8").Value = 92372
$ws.Range("F8").Value = 3878
$ws.Range("G8").Value = 218
$ws.Range("H8").Value = 24594

# Row 20 - Suiza
$ws.Range("E20").Value = 4551
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 1754

# Row 110 - Georgia
$ws.Range("E110").Value = 352
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 7

# Row 126 - Montenegro
$ws.Range("D126").Value = 233
$ws.Range("E126").Value = 82

# Row 171 - Siria
$ws.Range("B171").Value = 44
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 27
$ws.Range("E171").Value = 14
